$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "sCs"
$ws.Range("B2").Value = "Wnt5b"
$ws.Range("C2").Value = "Fzd8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8152283333333333
$ws.Range("H2").Value = 2.445685
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.027778333333333
$ws.Range("N2").Value = 3.083335
$ws.Range("O2").Value = 0.08020467841353289
$ws.Range("P2").Value = 0.08020467841353289
$ws.Range("Q2").Value = 0.8378740177194444
$ws.Range("R2").Value = 7.540866159475001
$ws.Range("S2").Value = 0.08020467841353289
$ws.Range("T2").Value = 0.08020467841353289

# Row 3
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Wnt5b"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8152283333333333
$ws.Range("H3").Value = 2.445685
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.273511666666667
$ws.Range("N3").Value = 21.820535
$ws.Range("O3").Value = 0.567602609669802
$ws.Range("P3").Value = 0.567602609669802
$ws.Range("Q3").Value = 5.929572793497222
$ws.Range("R3").Value = 53.366155141475
$ws.Range("S3").Value = 0.567602609669802
$ws.Range("T3").Value = 0.567602609669802

# Row 4
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Wnt5b"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8152283333333333
$ws.Range("H4").Value = 2.445685
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.513153666666667
$ws.Range("N4").Value = 13.539461
$ws.Range("O4").Value = 0.3521927119166651
$ws.Range("P4").Value = 0.3521927119166651
$ws.Range("Q4").Value = 3.679250741753889
$ws.Range("R4").Value = 33.113256675785
$ws.Range("S4").Value = 0.3521927119166651
$ws.Range("T4").Value = 0.3521927119166651

